# Scheduled market-data refresh: update computed Leve profit columns (H,I,J,K,L,M,N)
# across all item-crafting sheets. Values below are sourced from the latest
# Universalis price snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 24666.334
$ws.Range("J63").Value = 29499.5
$ws.Range("L63").Value = 29499.5
$ws.Range("N63").Value = -30747.5
$ws.Range("H66").Value = 24666.334
$ws.Range("J66").Value = 29499.5
$ws.Range("L66").Value = 88498.5
$ws.Range("N66").Value = -94738.5
$ws.Range("H98").Value = 3012.375
$ws.Range("I98").Value = 3700
$ws.Range("J98").Value = 2599.8
$ws.Range("K98").Value = 3700
$ws.Range("L98").Value = 2599.8
$ws.Range("M98").Value = -2202
$ws.Range("N98").Value = -5595.8
$ws.Range("H122").Value = 3012.375
$ws.Range("I122").Value = 3700
$ws.Range("J122").Value = 2599.8
$ws.Range("K122").Value = 11100
$ws.Range("L122").Value = 7799.400000000001
$ws.Range("M122").Value = -8650
$ws.Range("N122").Value = -12699.4
$ws.Range("H129").Value = 819
$ws.Range("I129").Value = 490
$ws.Range("J129").Value = 1202.8334
$ws.Range("K129").Value = 1470
$ws.Range("L129").Value = 3608.5002
$ws.Range("M129").Value = 3530
$ws.Range("N129").Value = -13608.5002
$ws.Range("H133").Value = 29093
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 29093
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 29093
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -39213

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 438
$ws.Range("I4").Value = 397.5
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 397.5
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -281.5
$ws.Range("N4").Value = -832
$ws.Range("H32").Value = 4258.46
$ws.Range("I32").Value = 3033.9302
$ws.Range("J32").Value = 11780.571
$ws.Range("K32").Value = 3033.9302
$ws.Range("L32").Value = 11780.571
$ws.Range("M32").Value = -2746.9302
$ws.Range("N32").Value = -12354.571

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1158.3182
$ws.Range("I20").Value = 1091.6
$ws.Range("J20").Value = 1301.2858
$ws.Range("K20").Value = 1091.6
$ws.Range("L20").Value = 1301.2858
$ws.Range("M20").Value = -844.5999999999999
$ws.Range("N20").Value = -1795.2858
$ws.Range("H98").Value = 29800
$ws.Range("J98").Value = 29800
$ws.Range("L98").Value = 29800
$ws.Range("N98").Value = -35790
$ws.Range("H105").Value = 2541.6216
$ws.Range("I105").Value = 2410.9678
$ws.Range("J105").Value = 3216.6667
$ws.Range("K105").Value = 2410.9678
$ws.Range("L105").Value = 3216.6667
$ws.Range("M105").Value = -663.9677999999999
$ws.Range("N105").Value = -6710.6667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71.333336
$ws.Range("I7").Value = 54
$ws.Range("J7").Value = 145
$ws.Range("K7").Value = 54
$ws.Range("L7").Value = 145
$ws.Range("M7").Value = 59
$ws.Range("N7").Value = -371
$ws.Range("H22").Value = 739.3077
$ws.Range("I22").Value = 460.66666
$ws.Range("J22").Value = 978.1429000000001
$ws.Range("K22").Value = 460.66666
$ws.Range("L22").Value = 978.1429000000001
$ws.Range("M22").Value = -110.66666
$ws.Range("N22").Value = -1678.1429
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H105").Value = 619.8
$ws.Range("I105").Value = 606.8182
$ws.Range("K105").Value = 606.8182
$ws.Range("M105").Value = 1140.1818

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 2960
$ws.Range("I124").Value = 1100
$ws.Range("J124").Value = 3425
$ws.Range("K124").Value = 3300
$ws.Range("L124").Value = 10275
$ws.Range("M124").Value = 1610
$ws.Range("N124").Value = -20095
$ws.Range("H132").Value = 3288.359
$ws.Range("J132").Value = 3948.2
$ws.Range("L132").Value = 35533.8
$ws.Range("N132").Value = -40593.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 40000
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40766
$ws.Range("H85").Value = 40000
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42652
$ws.Range("H116").Value = 39999
$ws.Range("J116").Value = 39999
$ws.Range("L116").Value = 39999
$ws.Range("N116").Value = -49177
$ws.Range("H126").Value = 4270.091
$ws.Range("I126").Value = 4268.5
$ws.Range("J126").Value = 4272
$ws.Range("K126").Value = 12805.5
$ws.Range("L126").Value = 12816
$ws.Range("M126").Value = -10335.5
$ws.Range("N126").Value = -17756

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2201.3809
$ws.Range("I7").Value = 1936
$ws.Range("J7").Value = 2732.1428
$ws.Range("K7").Value = 1936
$ws.Range("L7").Value = 2732.1428
$ws.Range("M7").Value = -1824
$ws.Range("N7").Value = -2956.1428
$ws.Range("H22").Value = 724
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405
$ws.Range("H27").Value = 724
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593
$ws.Range("H45").Value = 5231.4546
$ws.Range("I45").Value = 920
$ws.Range("J45").Value = 8824.333000000001
$ws.Range("K45").Value = 920
$ws.Range("L45").Value = 8824.333000000001
$ws.Range("M45").Value = -513
$ws.Range("N45").Value = -9638.333000000001
$ws.Range("H82").Value = 1447.25
$ws.Range("I82").Value = 1157.4706
$ws.Range("J82").Value = 1895.091
$ws.Range("K82").Value = 1157.4706
$ws.Range("L82").Value = 1895.091
$ws.Range("M82").Value = -796.4706000000001
$ws.Range("N82").Value = -2617.091
$ws.Range("H85").Value = 1447.25
$ws.Range("I85").Value = 1157.4706
$ws.Range("J85").Value = 1895.091
$ws.Range("K85").Value = 1157.4706
$ws.Range("L85").Value = 1895.091
$ws.Range("M85").Value = 90.5293999999999
$ws.Range("N85").Value = -4391.091
$ws.Range("H93").Value = 1295.5385
$ws.Range("I93").Value = 1291
$ws.Range("J93").Value = 1350
$ws.Range("K93").Value = 1291
$ws.Range("L93").Value = 1350
$ws.Range("M93").Value = -43
$ws.Range("N93").Value = -3846
$ws.Range("H100").Value = 43482116
$ws.Range("I100").Value = 5360.231
$ws.Range("J100").Value = 100001900
$ws.Range("K100").Value = 5360.231
$ws.Range("L100").Value = 100001900
$ws.Range("M100").Value = -4819.231
$ws.Range("N100").Value = -100002982
$ws.Range("H126").Value = 2201.3809
$ws.Range("I126").Value = 1936
$ws.Range("J126").Value = 2732.1428
$ws.Range("K126").Value = 5808
$ws.Range("L126").Value = 8196.428400000001
$ws.Range("M126").Value = -3338
$ws.Range("N126").Value = -13136.4284

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4046.9033
$ws.Range("I62").Value = 3840.5
$ws.Range("J62").Value = 4077.4814
$ws.Range("K62").Value = 3840.5
$ws.Range("L62").Value = 4077.4814
$ws.Range("M62").Value = -3216.5
$ws.Range("N62").Value = -5325.481400000001
$ws.Range("H65").Value = 4046.9033
$ws.Range("I65").Value = 3840.5
$ws.Range("J65").Value = 4077.4814
$ws.Range("K65").Value = 19202.5
$ws.Range("L65").Value = 20387.407
$ws.Range("M65").Value = -16082.5
$ws.Range("N65").Value = -26627.407
$ws.Range("H132").Value = 22601.06
$ws.Range("I132").Value = 34609.332
$ws.Range("J132").Value = 3640.6316
$ws.Range("K132").Value = 103827.996
$ws.Range("L132").Value = 10921.8948
$ws.Range("M132").Value = -101297.996
$ws.Range("N132").Value = -15981.8948

